$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number 32 -> 33, week range 8/7-8/13/2023 -> 8/14-8/20/2023 ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Crime-statistics table updates (rows 14-29) ---
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L14").Value = 200
$ws.Range("N14").Value = -50
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = -38.095238095238
$ws.Range("I16").Value = 102
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = -15.702479338843
$ws.Range("L16").Value = 82.142857142857
$ws.Range("M16").Value = -2.857142857142
$ws.Range("N16").Value = -72.727272727272
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -34.782608695652
$ws.Range("J17").Value = 142
$ws.Range("K17").Value = 0.704225352112
$ws.Range("L17").Value = 47.422680412371
$ws.Range("N17").Value = -20.111731843575
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 98
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = 11.363636363636
$ws.Range("L18").Value = 84.905660377358
$ws.Range("M18").Value = 10.112359550561
$ws.Range("N18").Value = -86.274509803921
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -36.842105263157
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -24.590163934426
$ws.Range("I19").Value = 413
$ws.Range("J19").Value = 436
$ws.Range("K19").Value = -5.275229357798
$ws.Range("L19").Value = 51.838235294117
$ws.Range("M19").Value = 84.375
$ws.Range("N19").Value = 49.637681159420
$ws.Range("C20").Value = 8
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 277.777777777778
$ws.Range("I20").Value = 219
$ws.Range("K20").Value = 32.727272727272
$ws.Range("L20").Value = 208.450704225352
$ws.Range("M20").Value = 173.75
$ws.Range("N20").Value = -81.345826235093
$ws.Range("D21").Value = 29
$ws.Range("E21").Value = -6.896551724137
$ws.Range("F21").Value = 114
$ws.Range("H21").Value = -9.523809523809
$ws.Range("I21").Value = 984
$ws.Range("J21").Value = 967
$ws.Range("K21").Value = 1.758014477766
$ws.Range("L21").Value = 75.088967971530
$ws.Range("M21").Value = 66.216216216216
$ws.Range("N21").Value = -64.008778346744
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 57
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 14
$ws.Range("L23").Value = 137.5
$ws.Range("M23").Value = 103.571428571429
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 29.166666666666
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -4.629629629629
$ws.Range("I24").Value = 703
$ws.Range("J24").Value = 764
$ws.Range("K24").Value = -7.984293193717
$ws.Range("L24").Value = 29.465930018416
$ws.Range("M24").Value = 38.385826771653
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 12.903225806451
$ws.Range("I25").Value = 207
$ws.Range("J25").Value = 221
$ws.Range("K25").Value = -6.334841628959
$ws.Range("L25").Value = 26.219512195122
$ws.Range("M25").Value = -16.867469879518
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 2
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 17
$ws.Range("K26").Value = -5.555555555555
$ws.Range("L26").Value = 13.333333333333
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 19
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -20.833333333333
$ws.Range("L27").Value = 5.555555555555
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("F28").NumberFormat = "General"
$ws.Range("H28").Value = -100
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -20
$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("F29").NumberFormat = "General"
$ws.Range("H29").Value = -100
$ws.Range("L29").Value = -42.857142857142
$ws.Range("M29").Value = -20
